$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 86 (shifts existing rows 86-93 down to 87-94,
# and shifts their formatting / merged cells along with them).
$ws.Rows("86:86").Insert()

# The row that used to be row 86 is now row 87, with the correct cell
# styles already in place. Copy those styles onto the freshly inserted
# (blank/default-styled) row 86 so the new row matches the sheet's table
# formatting (merges, borders, number formats, etc.).
$ws.Range("A87:Q87").Copy()
$ws.Range("A86:Q86").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row with the new product line: "قصافات كبار".
$ws.Range("A86").Value = 80
$ws.Range("C86").Value = "قصافات كبار"
$ws.Range("H86").Value = "17:0"
$ws.Range("L86").Value = "0"
$ws.Range("N86").Value = "25.00"
$ws.Range("P86").Value = "25.0000"
$ws.Range("Q86").Value = "1:0"

# Update the running total (old total 7639.655 + new row's 25.00 = 7664.655).
# This cell is now on row 93 after the insert shifted it down from row 92.
$ws.Range("P93").Value = 7664.6549999999997

# Update the generated timestamp in the report footer (now on row 94).
$ws.Range("A94").Value = "Thursday, 18 September, 2025 10:56 PM"
